$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column F ("język"), shifting columns G:J left to F:I
$ws.Range("F:F").Delete()
